$d = $word.ActiveDocument

# 1. "Motor 2 (up" + "): (17,23) = in java (0,4)" -> merge into a single run's text
$d.Content.Find.Execute("Motor 2 (up): (17,23) = in java (0,4)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Motor 2 (up): (17,23) = in java (0,4)", 2)

# 2. "Motor 3 (right" + "): (9,7) = in java (13,11)" -> merge into a single run's text
$d.Content.Find.Execute("Motor 3 (right): (9,7) = in java (13,11)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Motor 3 (right): (9,7) = in java (13,11)", 2)

# 3. Append a new paragraph after the last one ("pwm (up): 18 = in java 1"),
#    with the same spell-check-marker (proofErr) structure Word would add
#    around the non-dictionary tokens "mosi,sclk" and "pinnen".
$range = $d.Content
$range.Collapse(0)
$xml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:r><w:t>(10,11) = in java(</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>mosi,sclk</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> of 12,14 == de </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>pinnen</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> van motor 4!!)</w:t></w:r>
</w:p>
"@
$range.InsertXML($xml)
